# Refresh the coin Price (D) / Volume(1h) (E) columns with the latest scrape,
# and fix the WrappedEther/Polygon row ordering (rows 13 & 14 swapped places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.646.89'
$ws.Range('E2').Value = '  -2.14%  '
$ws.Range('D3').Value = '1.808.63'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = "'231.68"
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').Value = "'0.603"
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').Value = "'39.31"
$ws.Range('E8').Value = '  -9.06%  '
$ws.Range('D9').Value = "'0.327"
$ws.Range('E9').Value = '  +6.10%  '
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('D12').Value = '2.071.89'
$ws.Range('E12').Value = '  -1.68%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = "'0.673"
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.806.12'
$ws.Range('E14').Value = '  -1.94%  '
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('E16').Value = '  -1.87%  '
$ws.Range('D17').Value = '34.641.75'
$ws.Range('D18').Value = "'69.39"
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('D19').Value = '0.0₃0784'
$ws.Range('E19').Value = '  -1.77%  '
$ws.Range('D20').Value = "'240.21"
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('D21').Value = "'11.88"
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('D22').Value = "'4.70"
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('D24').Value = "'2.24"
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('D25').Value = "'171.52"
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('D26').Value = "'7.70"
$ws.Range('E26').Value = '  -2.34%  '
$ws.Range('D27').Value = "'17.18"
$ws.Range('E27').Value = '  -3.02%  '
$ws.Range('D28').Value = "'0.120"
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('E29').Value = '  -1.04%  '
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').Value = "'4.02"
$ws.Range('E31').Value = '  +2.58%  '
$ws.Range('D32').Value = "'0.0545"
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').Value = "'3.94"
$ws.Range('E33').Value = '  -2.63%  '
$ws.Range('D34').Value = "'1.32"
$ws.Range('E34').Value = '  +21.78%  '
$ws.Range('E35').Value = '  -3.59%  '
$ws.Range('D36').Value = "'0.698"
$ws.Range('E36').Value = '  +1.89%  '
$ws.Range('D37').Value = "'91.35"
$ws.Range('E37').Value = '  -2.90%  '
$ws.Range('D38').Value = "'1.32"
$ws.Range('E38').Value = '  +4.38%  '
$ws.Range('D39').Value = '1.327.39'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('E42').Value = '  -4.44%  '
$ws.Range('D43').Value = "'14.16"
$ws.Range('E43').Value = '  -6.87%  '
$ws.Range('E44').Value = '  -9.38%  '
$ws.Range('E45').Value = '  -4.86%  '
$ws.Range('E46').Value = '  +1.04%  '
$ws.Range('E47').Value = '  -1.24%  '
$ws.Range('D48').Value = '1.999.18'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('E50').Value = '  +7.30%  '
$ws.Range('D51').Value = "'98.07"
$ws.Range('E51').Value = '  -4.89%  '
